# Apply updated Model_Home_win and Model_home_win_probability values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 2; C = $null; D = 0.4918893277645111 },
    @{ Row = 3; C = 1; D = 0.6692629456520081 },
    @{ Row = 4; C = $null; D = 0.4968279898166656 },
    @{ Row = 5; C = 0; D = 0.3787283301353455 },
    @{ Row = 6; C = $null; D = 0.5030312538146973 },
    @{ Row = 7; C = $null; D = 0.5558325052261353 },
    @{ Row = 8; C = $null; D = 0.4457062184810638 },
    @{ Row = 9; C = 1; D = 0.5066936612129211 },
    @{ Row = 10; C = $null; D = 0.3962562084197998 },
    @{ Row = 11; C = $null; D = 0.5024792551994324 },
    @{ Row = 12; C = $null; D = 0.5051472783088684 },
    @{ Row = 13; C = 1; D = 0.5567144155502319 },
    @{ Row = 14; C = $null; D = 0.5188864469528198 },
    @{ Row = 15; C = 1; D = 0.9265148639678955 },
    @{ Row = 16; C = 1; D = 0.7495768666267395 },
    @{ Row = 17; C = $null; D = 0.3702545762062073 },
    @{ Row = 18; C = $null; D = 0.5541085004806519 },
    @{ Row = 19; C = 1; D = 0.5323666334152222 },
    @{ Row = 20; C = $null; D = 0.5690140128135681 },
    @{ Row = 21; C = 0; D = 0.4938535988330841 },
    @{ Row = 22; C = $null; D = 0.4827529788017273 },
    @{ Row = 23; C = $null; D = 0.4715619385242462 },
    @{ Row = 24; C = $null; D = 0.5267050266265869 },
    @{ Row = 25; C = $null; D = 0.5497201085090637 },
    @{ Row = 26; C = $null; D = 0.5774034857749939 },
    @{ Row = 27; C = $null; D = 0.4997602999210358 },
    @{ Row = 28; C = 0; D = 0.4899245500564575 },
    @{ Row = 29; C = 1; D = 0.5293914675712585 },
    @{ Row = 30; C = $null; D = 0.4967174530029297 },
    @{ Row = 31; C = 1; D = 0.5032103061676025 },
    @{ Row = 32; C = $null; D = 0.4749829173088074 },
    @{ Row = 33; C = $null; D = 0.4550463855266571 },
    @{ Row = 34; C = 1; D = 0.5627449750900269 },
    @{ Row = 35; C = 0; D = 0.482975572347641 },
    @{ Row = 36; C = $null; D = 0.4746008813381195 },
    @{ Row = 37; C = 0; D = 0.2514096200466156 },
    @{ Row = 38; C = 0; D = 0.4619325995445251 },
    @{ Row = 39; C = $null; D = 0.6663196086883545 },
    @{ Row = 40; C = 0; D = 0.4759460985660553 },
    @{ Row = 41; C = $null; D = 0.470602810382843 },
    @{ Row = 42; C = $null; D = 0.5046419501304626 },
    @{ Row = 43; C = $null; D = 0.4914324581623077 },
    @{ Row = 44; C = 0; D = 0.4040604829788208 },
    @{ Row = 45; C = $null; D = 0.4991129338741302 },
    @{ Row = 46; C = 1; D = 0.5260767340660095 },
    @{ Row = 47; C = 1; D = 0.5174676775932312 },
    @{ Row = 48; C = $null; D = 0.5244544148445129 },
    @{ Row = 49; C = 0; D = 0.3714616298675537 }
)

foreach ($u in $updates) {
    if ($null -ne $u.C) {
        $ws.Cells.Item($u.Row, 3).Value = $u.C
    }
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
